$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5779385566711426
$ws.Range("B1").Value = 1.184473872184753
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.735146760940552
$ws.Range("E1").Value = 1.453469395637512
